# Updated cryptos list values (price + 1h volume change columns),
# matching the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.150.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.076.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5225"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4343"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "55.05"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09296"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.165"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.401"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.819"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.051.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "100.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001151"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06696"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.247"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.189.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.308"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.758"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.455"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.116"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1040"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.622"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.208"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.894"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02586"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.761"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06665"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6871"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.325"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2176"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6659"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.329"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.311"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.617"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000349"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.201"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07199"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
